$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) The "_GoBack" bookmark moves from the "Header Comment, and Formatting"
#    heading down into the new "Runtime complexity analysis" heading.  Word
#    only ever keeps a single "_GoBack" bookmark, so remove the stale one now
#    (it gets re-created at the new location once that text exists below).
# ---------------------------------------------------------------------------
$ErrorActionPreference = "SilentlyContinue"
$goBack = $d.Bookmarks.Item("_GoBack")
if ($goBack) {
    $goBack.Delete()
}
$ErrorActionPreference = "Continue"

# ---------------------------------------------------------------------------
# 2) Replace the placeholder "To be completed" paragraph text with the real
#    submission-deadline sentence.
# ---------------------------------------------------------------------------
$deliverP = $d.Paragraphs.Item($d.Paragraphs.Count)
$deliverP.Range.Text = "The final version of your program must be uploaded through Canvas no later than midnight on "

$tail = $deliverP.Range
$tail.Collapse(0)
$tail.InsertAfter("Wednesday, February 13, 2019")

$tail2 = $deliverP.Range
$tail2.Collapse(0)
$tail2.InsertAfter(".  ")

# ---------------------------------------------------------------------------
# 3) Append the new "Grading Criteria" section after it.
# ---------------------------------------------------------------------------
$cur = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$cur.InsertParagraphAfter()
$gradingCriteria = $d.Paragraphs.Item($d.Paragraphs.Count)
$gradingCriteria.Range.ParagraphFormat.Style = "Heading2"
$gradingCriteria.Range.Text = "Grading Criteria"

$cur = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$cur.InsertParagraphAfter()
$judged = $d.Paragraphs.Item($d.Paragraphs.Count)
$judged.Range.Text = "Your assignment will be judged by the following criteria:"

# Test Cases heading (80pts)
$cur = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$cur.InsertParagraphAfter()
$testCasesHeading = $d.Paragraphs.Item($d.Paragraphs.Count)
$testCasesHeading.Range.ParagraphFormat.Style = "Heading3"
$testCasesHeading.Range.Text = "Test Cases (80pts)"

# Test Cases list item
$cur = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$cur.InsertParagraphAfter()
$testCasesItem = $d.Paragraphs.Item($d.Paragraphs.Count)
$testCasesItem.Range.ParagraphFormat.Style = "ListParagraph"
$testCasesItem.Range.Text = "Your program successfully passes all test cases"
$testCasesItem.Range.ListFormat.ApplyNumberedDefault()
$testCasesItem.Range.ParagraphFormat.SpaceAfter = 10
$testCasesItem.Range.ParagraphFormat.LineSpacingRule = 5
$testCasesItem.Range.ParagraphFormat.LineSpacing = 13.8
$testCasesItem.Range.ParagraphFormat.Alignment = 0

# Reflection essay heading (10pts)
$cur = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$cur.InsertParagraphAfter()
$reflectionHeading = $d.Paragraphs.Item($d.Paragraphs.Count)
$reflectionHeading.Range.ParagraphFormat.Style = "Heading3"
$reflectionHeading.Range.Text = "Reflection essay (10pts)"

# Reflection essay list item
$cur = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$cur.InsertParagraphAfter()
$reflectionItem = $d.Paragraphs.Item($d.Paragraphs.Count)
$reflectionItem.Range.ParagraphFormat.Style = "ListParagraph"
$reflectionItem.Range.Text = "Your reflection meets the minimum requirements as specified earlier in this document. "
$reflectionItem.Range.ListFormat.ApplyNumberedDefault()

# Runtime complexity analysis heading (10pts) - also hosts the _GoBack bookmark
$cur = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$cur.InsertParagraphAfter()
$runtimeHeading = $d.Paragraphs.Item($d.Paragraphs.Count)
$runtimeHeading.Range.ParagraphFormat.Style = "Heading3"
$runtimeHeading.Range.Text = "Runtime complexity analysis (10pts)"

# Runtime complexity analysis list item
$cur = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$cur.InsertParagraphAfter()
$runtimeItem = $d.Paragraphs.Item($d.Paragraphs.Count)
$runtimeItem.Range.ParagraphFormat.Style = "ListParagraph"
$runtimeItem.Range.Text = "You correctly identify the runtime complexity of your functions."
$runtimeItem.Range.ListFormat.ApplyNumberedDefault()

# Final, empty Heading3 paragraph (becomes the new last paragraph of the body).
$cur = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$cur.InsertParagraphAfter()
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.ParagraphFormat.Style = "Heading3"

# ---------------------------------------------------------------------------
# 4) Re-create the "_GoBack" bookmark at its new home: right after the "10"
#    in "Runtime complexity analysis (10pts)", matching the source edit.
# ---------------------------------------------------------------------------
$markerPos = $runtimeHeading.Range.Start + "Runtime complexity analysis (10".Length
$markerRange = $d.Range($markerPos, $markerPos)
$d.Bookmarks.Add("_GoBack", $markerRange)

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
